$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

$styledRange = $excel.Union($ws.Range("B1"), $ws.Range("A2"))
$styledRange.Font.Bold = $true
$styledRange.HorizontalAlignment = -4108
$styledRange.VerticalAlignment = -4160
$styledRange.Borders.LineStyle = 1
$styledRange.Borders.Weight = 2
